$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.01570068801409834
$ws.Range("C2").Value = 2.201232587161126
$ws.Range("D2").Value = 18.29599404310352
$ws.Range("E2").Value = 4.277381680783645
$ws.Range("F2").Value = 4.382982421602256

# Row 3
$ws.Range("B3").Value = 0.9935082891590984
$ws.Range("C3").Value = 2.320399451120858
$ws.Range("D3").Value = 20.93761159639632
$ws.Range("E3").Value = 4.575763498739453
$ws.Range("F3").Value = 4.571715127081522

# Row 4
$ws.Range("B4").Value = -0.1734791790674897
$ws.Range("C4").Value = 0.9918456812746188
$ws.Range("D4").Value = 3.811757031605553
$ws.Range("E4").Value = 1.952372154996468
$ws.Range("F4").Value = 1.997937521471647
$ws.Range("G4").Value = 19
